$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0
$ws.Range("A4").Value = 0
$ws.Range("A5").Value = 0
$ws.Range("A6").Value = 0
$ws.Range("A7").Value = 0
$ws.Range("A8").Value = 0
